$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = 0.5088967971530249
$summary.Range("C2").Value = 0.07823129251700681
$summary.Range("D2").Value = 0.8214285714285714
$summary.Range("E2").Value = 0.1428571428571428
$summary.Range("F2").Value = 0.2832512315270936
$summary.Range("G2").Value = 0.6016096579476862
$summary.Range("H2").Value = 0.7286650615302301
$summary.Range("I2").Value = 23
$summary.Range("J2").Value = 271
$summary.Range("K2").Value = 263
$summary.Range("L2").Value = 5

# --- Sheet: Classification Report ---
$cr = $wb.Worksheets.Item("Classification Report")

# Row 2 ("0")
$cr.Range("B2").Value = 0.9813432835820896
$cr.Range("C2").Value = 0.4925093632958801
$cr.Range("D2").Value = 0.6558603491271821

# Row 3 ("1")
$cr.Range("B3").Value = 0.07823129251700681
$cr.Range("C3").Value = 0.8214285714285714
$cr.Range("D3").Value = 0.1428571428571428

# Row 4 ("accuracy")
$cr.Range("B4").Value = 0.5088967971530249
$cr.Range("C4").Value = 0.5088967971530249
$cr.Range("D4").Value = 0.5088967971530249
$cr.Range("E4").Value = 0.5088967971530249

# Row 5 ("macro avg")
$cr.Range("B5").Value = 0.5297872880495482
$cr.Range("C5").Value = 0.6569689673622258
$cr.Range("D5").Value = 0.3993587459921625

# Row 6 ("weighted avg")
$cr.Range("B6").Value = 0.9363483801126548
$cr.Range("C6").Value = 0.5088967971530249
$cr.Range("D6").Value = 0.6303014705229808

# --- Sheet: Confusion Matrix ---
$cm = $wb.Worksheets.Item("Confusion Matrix")

# Row 2 ("Actual 0")
$cm.Range("B2").Value = 263
$cm.Range("C2").Value = 271

# Row 3 ("Actual 1")
$cm.Range("B3").Value = 5
$cm.Range("C3").Value = 23
